$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 9).Value = "24-10-2025 00:00:00"
$ws.Cells.Item(1, 11).Value = "24-10-2025 00:00:00"
$ws.Cells.Item(1, 13).Value = 45954
$ws.Cells.Item(90, 6).Value = 123
$ws.Cells.Item(90, 7).Value = 8215.17
$ws.Cells.Item(98, 6).Value = 154
$ws.Cells.Item(98, 7).Value = 20782.3
$ws.Cells.Item(124, 2).Value = 368979.97
$ws.Cells.Item(218, 6).Value = 167
$ws.Cells.Item(218, 7).Value = 14508.96
$ws.Cells.Item(225, 2).Value = 71025.67999999999
$ws.Cells.Item(228, 6).Value = 9
$ws.Cells.Item(228, 7).Value = 166.5
$ws.Cells.Item(232, 2).Value = 2126.63
$ws.Cells.Item(268, 6).Value = 55
$ws.Cells.Item(268, 7).Value = 2748.35
$ws.Cells.Item(310, 2).Value = 148409.58
$ws.Cells.Item(313, 2).Value = 57854
$ws.Cells.Item(313, 6).Value = 2
$ws.Cells.Item(313, 7).Value = 611.6799999999999
$ws.Cells.Item(314, 2).Value = 62997
$ws.Cells.Item(314, 6).Value = 0
$ws.Cells.Item(314, 7).Value = 0
$ws.Cells.Item(339, 6).Value = 5
$ws.Cells.Item(339, 7).Value = 603.25
$ws.Cells.Item(351, 2).Value = 57802
$ws.Cells.Item(351, 5).Value = 162.71
$ws.Cells.Item(351, 6).Value = -79
$ws.Cells.Item(351, 7).Value = -11334.92
$ws.Cells.Item(352, 2).Value = 63531
$ws.Cells.Item(352, 5).Value = 152.53
$ws.Cells.Item(352, 6).Value = 80
$ws.Cells.Item(352, 7).Value = 11478.4
$ws.Cells.Item(355, 2).Value = 63510
$ws.Cells.Item(355, 5).Value = 50.66
$ws.Cells.Item(355, 6).Value = 154
$ws.Cells.Item(355, 7).Value = 7336.56
$ws.Cells.Item(356, 2).Value = 55356
$ws.Cells.Item(356, 5).Value = 54.04
$ws.Cells.Item(356, 6).Value = -158
$ws.Cells.Item(356, 7).Value = -7527.12
$ws.Cells.Item(372, 2).Value = 57885
$ws.Cells.Item(372, 5).Value = 62.28
$ws.Cells.Item(372, 6).Value = 4
$ws.Cells.Item(372, 7).Value = 208.52
$ws.Cells.Item(373, 2).Value = 63652
$ws.Cells.Item(373, 5).Value = 55.42
$ws.Cells.Item(373, 6).Value = 198
$ws.Cells.Item(373, 7).Value = 10321.74
$ws.Cells.Item(375, 2).Value = 63563
$ws.Cells.Item(375, 5).Value = 119.04
$ws.Cells.Item(375, 6).Value = 2
$ws.Cells.Item(375, 7).Value = 223.92
$ws.Cells.Item(376, 2).Value = 61605
$ws.Cells.Item(376, 5).Value = 133.78
$ws.Cells.Item(376, 6).Value = -13
$ws.Cells.Item(376, 7).Value = -1455.48
$ws.Cells.Item(382, 2).Value = 63560
$ws.Cells.Item(382, 5).Value = 134.87
$ws.Cells.Item(382, 6).Value = 22
$ws.Cells.Item(382, 7).Value = 2790.92
$ws.Cells.Item(383, 2).Value = 60325
$ws.Cells.Item(383, 5).Value = 151.57
$ws.Cells.Item(383, 6).Value = -102
$ws.Cells.Item(383, 7).Value = -12939.72
$ws.Cells.Item(389, 2).Value = 57817
$ws.Cells.Item(389, 6).Value = 3
$ws.Cells.Item(389, 7).Value = 239.43
$ws.Cells.Item(390, 2).Value = 62865
$ws.Cells.Item(390, 6).Value = 46
$ws.Cells.Item(390, 7).Value = 3671.26
$ws.Cells.Item(408, 6).Value = 2426
$ws.Cells.Item(408, 7).Value = 51018.78
$ws.Cells.Item(411, 6).Value = 85
$ws.Cells.Item(411, 7).Value = 13685
$ws.Cells.Item(419, 2).Value = 57856
$ws.Cells.Item(419, 6).Value = 2
$ws.Cells.Item(419, 7).Value = 342.66
$ws.Cells.Item(420, 2).Value = 63007
$ws.Cells.Item(420, 6).Value = 858
$ws.Cells.Item(420, 7).Value = 147001.14
$ws.Cells.Item(421, 2).Value = 57857
$ws.Cells.Item(421, 6).Value = 3
$ws.Cells.Item(421, 7).Value = 453.51
$ws.Cells.Item(422, 2).Value = 63008
$ws.Cells.Item(422, 6).Value = 444
$ws.Cells.Item(422, 7).Value = 67119.48
$ws.Cells.Item(431, 2).Value = 63102
$ws.Cells.Item(431, 3).Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Cells.Item(431, 6).Value = 4
$ws.Cells.Item(431, 7).Value = 237.88
$ws.Cells.Item(432, 2).Value = 53082
$ws.Cells.Item(432, 3).Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Cells.Item(432, 6).Value = 1
$ws.Cells.Item(432, 7).Value = 59.47
$ws.Cells.Item(434, 2).Value = 540785.99
$ws.Cells.Item(579, 2).Value = 53757
$ws.Cells.Item(579, 5).Value = 16.08
$ws.Cells.Item(579, 6).Value = -159
$ws.Cells.Item(579, 7).Value = -2138.55
$ws.Cells.Item(580, 2).Value = 65069
$ws.Cells.Item(580, 5).Value = 14.3
$ws.Cells.Item(580, 6).Value = 85
$ws.Cells.Item(580, 7).Value = 1143.25
$ws.Cells.Item(583, 2).Value = 53263
$ws.Cells.Item(583, 5).Value = 15.29
$ws.Cells.Item(583, 6).Value = -309
$ws.Cells.Item(583, 7).Value = -3958.29
$ws.Cells.Item(584, 2).Value = 65066
$ws.Cells.Item(584, 5).Value = 13.61
$ws.Cells.Item(584, 6).Value = 249
$ws.Cells.Item(584, 7).Value = 3189.69
$ws.Cells.Item(586, 2).Value = 64915
$ws.Cells.Item(586, 5).Value = 20.98
$ws.Cells.Item(586, 6).Value = 13
$ws.Cells.Item(586, 7).Value = 256.49
$ws.Cells.Item(587, 2).Value = 45695
$ws.Cells.Item(587, 5).Value = 23.58
$ws.Cells.Item(587, 6).Value = -36
$ws.Cells.Item(587, 7).Value = -710.28
$ws.Cells.Item(604, 2).Value = 53595
$ws.Cells.Item(604, 5).Value = 17.61
$ws.Cells.Item(604, 6).Value = -335
$ws.Cells.Item(604, 7).Value = -4934.55
$ws.Cells.Item(605, 2).Value = 65067
$ws.Cells.Item(605, 5).Value = 15.65
$ws.Cells.Item(605, 6).Value = 311
$ws.Cells.Item(605, 7).Value = 4581.03
$ws.Cells.Item(649, 6).Value = 36
$ws.Cells.Item(649, 7).Value = 3789.72
$ws.Cells.Item(651, 2).Value = 42104.07
$ws.Cells.Item(854, 6).Value = 324
$ws.Cells.Item(854, 7).Value = 46374.12
$ws.Cells.Item(857, 6).Value = 150
$ws.Cells.Item(857, 7).Value = 19965
$ws.Cells.Item(867, 6).Value = 782
$ws.Cells.Item(867, 7).Value = 105577.82
$ws.Cells.Item(869, 6).Value = 766
$ws.Cells.Item(869, 7).Value = 92463.86
$ws.Cells.Item(870, 2).Value = 353754.38
$ws.Cells.Item(889, 2).Value = 65079
$ws.Cells.Item(889, 6).Value = 21
$ws.Cells.Item(889, 7).Value = 858.27
$ws.Cells.Item(890, 2).Value = 65362
$ws.Cells.Item(890, 6).Value = 88
$ws.Cells.Item(890, 7).Value = 3596.56
$ws.Cells.Item(962, 2).Value = 4755113.06
$ws.Cells.Item(963, 2).Value = 4755113.06
